# Insert a new weekly record at row 200 (Fruta, Feria Lagunitas de Puerto Montt - Piña).
# Inserting the row shifts every existing row from 200 downward (old row 200 -> 201,
# ..., old row 241 -> 242), which matches the diff's "everything moved down by one" shape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(200).Insert()

$ws.Cells.Item(200, 1).Value  = 4
$ws.Cells.Item(200, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(200, 3).Value  = "Los Lagos"
$ws.Cells.Item(200, 4).Value  = 44694
$ws.Cells.Item(200, 5).Value  = 10
$ws.Cells.Item(200, 6).Value  = "Fruta"
$ws.Cells.Item(200, 7).Value  = 100108
$ws.Cells.Item(200, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(200, 9).Value  = 100108005
$ws.Cells.Item(200, 10).Value = "Piña"
$ws.Cells.Item(200, 11).Value = "Caramelo"
$ws.Cells.Item(200, 12).Value = "Tercera"
$ws.Cells.Item(200, 13).Value = 400
$ws.Cells.Item(200, 14).Value = 22000
$ws.Cells.Item(200, 15).Value = 23000
$ws.Cells.Item(200, 16).Value = 22500
$ws.Cells.Item(200, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(200, 18).Value = "Ecuador"
$ws.Cells.Item(200, 19).Value = 1406
$ws.Cells.Item(200, 20).Value = 16
